# Fix Training Data Issue (#48)
# The "Date" column (BF) was populated with the source filename-derived
# string "6-11-2011-12" for every team row instead of the actual game
# date. Correct it to the real ISO date "2012-06-11".
#
# The replacement text looks like a date, so a plain Value assignment
# would get auto-parsed by Excel into a date serial number. Prefixing
# with an apostrophe forces it to stay literal text (as it was before),
# matching the original inline-string "Date" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDate = "6-11-2011-12"
$newDate = "2012-06-11"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Range("BF$row")
    if ($cell.Text -eq $oldDate) {
        $cell.Value = "'" + $newDate
    }
}
